# Update the "想去人数" (want-to-go count) values in column F for the
# sheets that carry the event data ("展览" and "全部类型"). Both sheets
# contain identical row layouts, so the same updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1042
    3  = 727
    5  = 21
    6  = 1091
    8  = 1648
    9  = 6057
    10 = 474
    11 = 341
    12 = 273
    13 = 82
    14 = 356
    15 = 127
    16 = 4786
    17 = 256
    18 = 1259
    19 = 132
    20 = 107
    21 = 217
    22 = 96
    23 = 249
    24 = 94
    26 = 3
    27 = 90
    28 = 374
    29 = 67
    32 = 37
    34 = 19
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
